# Auto-generated edit script: applies numeric updates to the Goblin Profits workbook
# across all 8 profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Update cell values on ALC
$ws.Range("H32").Value = 7598.4
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 8998
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 8998
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -9650
$ws.Range("H33").Value = 702.3
$ws.Range("I33").Value = 163.35715
$ws.Range("K33").Value = 163.35715
$ws.Range("M33").Value = 65.64285000000001
$ws.Range("H38").Value = 8968.666999999999
$ws.Range("I38").Value = 11048.546
$ws.Range("K38").Value = 33145.638
$ws.Range("M38").Value = -32773.638
$ws.Range("H40").Value = 2789.8
$ws.Range("J40").Value = 3342.5715
$ws.Range("L40").Value = 3342.5715
$ws.Range("N40").Value = -3692.5715
$ws.Range("H64").Value = 4985.5645
$ws.Range("I64").Value = 3415.2827
$ws.Range("J64").Value = 9500.125
$ws.Range("K64").Value = 3415.2827
$ws.Range("L64").Value = 9500.125
$ws.Range("M64").Value = -3167.2827
$ws.Range("N64").Value = -9996.125
$ws.Range("H67").Value = 4985.5645
$ws.Range("I67").Value = 3415.2827
$ws.Range("J67").Value = 9500.125
$ws.Range("K67").Value = 3415.2827
$ws.Range("L67").Value = 9500.125
$ws.Range("M67").Value = -2557.2827
$ws.Range("N67").Value = -11216.125
$ws.Range("H69").Value = 5875
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("H72").Value = 5875
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("H74").Value = 2999.8
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 2999.8
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H88").Value = 4672.409
$ws.Range("I88").Value = 1575
$ws.Range("J88").Value = 4982.15
$ws.Range("K88").Value = 1575
$ws.Range("L88").Value = 4982.15
$ws.Range("M88").Value = -1169
$ws.Range("N88").Value = -5794.15
$ws.Range("H91").Value = 4672.409
$ws.Range("I91").Value = 1575
$ws.Range("J91").Value = 4982.15
$ws.Range("K91").Value = 1575
$ws.Range("L91").Value = 4982.15
$ws.Range("M91").Value = -171
$ws.Range("N91").Value = -7790.15
$ws.Range("H94").Value = 2179
$ws.Range("I94").Value = 2179
$ws.Range("K94").Value = 2179
$ws.Range("M94").Value = -1728
$ws.Range("H125").Value = 2269.2727
$ws.Range("I125").Value = 1044.6666
$ws.Range("J125").Value = 3738.8
$ws.Range("K125").Value = 9401.999400000001
$ws.Range("L125").Value = 33649.2
$ws.Range("M125").Value = -6941.999400000001
$ws.Range("N125").Value = -38569.2
$ws.Range("H132").Value = 2224270
$ws.Range("I132").Value = 2012.5405
$ws.Range("K132").Value = 6037.6215
$ws.Range("M132").Value = -3507.6215

# Clear cells that no longer hold a value on ALC
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Update cell values on ARM
$ws.Range("H2").Value = 2130.182
$ws.Range("I2").Value = 490.2857
$ws.Range("K2").Value = 490.2857
$ws.Range("M2").Value = -377.2857
$ws.Range("H61").Value = 786950.0600000001
$ws.Range("I61").Value = 5657.2856
$ws.Range("J61").Value = 6255999.5
$ws.Range("K61").Value = 5657.2856
$ws.Range("L61").Value = 6255999.5
$ws.Range("M61").Value = -5445.2856
$ws.Range("N61").Value = -6256423.5
$ws.Range("H88").Value = 2993.65
$ws.Range("J88").Value = 2897.6155
$ws.Range("L88").Value = 2897.6155
$ws.Range("N88").Value = -3709.6155
$ws.Range("H91").Value = 2993.65
$ws.Range("J91").Value = 2897.6155
$ws.Range("L91").Value = 2897.6155
$ws.Range("N91").Value = -5705.6155
$ws.Range("H92").Value = 33869.7
$ws.Range("J92").Value = 34310.777
$ws.Range("L92").Value = 34310.777
$ws.Range("N92").Value = -39302.777
$ws.Range("H97").Value = 417.77777
$ws.Range("I97").Value = 417.77777
$ws.Range("K97").Value = 417.77777
$ws.Range("M97").Value = 78.22223000000002
$ws.Range("H116").Value = 2130.182
$ws.Range("I116").Value = 490.2857
$ws.Range("K116").Value = 490.2857
$ws.Range("M116").Value = 1803.7143
$ws.Range("H122").Value = 6185557
$ws.Range("I122").Value = 10111230
$ws.Range("J122").Value = 16642.285
$ws.Range("K122").Value = 30333690
$ws.Range("L122").Value = 49926.855
$ws.Range("M122").Value = -30331240
$ws.Range("N122").Value = -54826.855
$ws.Range("H132").Value = 5514.24
$ws.Range("I132").Value = 6326.476
$ws.Range("K132").Value = 18979.428
$ws.Range("M132").Value = -16449.428
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H136").Value = 786950.0600000001
$ws.Range("I136").Value = 5657.2856
$ws.Range("J136").Value = 6255999.5
$ws.Range("K136").Value = 16971.8568
$ws.Range("L136").Value = 18767998.5
$ws.Range("M136").Value = -14421.8568
$ws.Range("N136").Value = -18773098.5

# Clear cells that no longer hold a value on ARM
$ws.Range("N134").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Update cell values on BSM
$ws.Range("H3").Value = 2130.182
$ws.Range("I3").Value = 490.2857
$ws.Range("K3").Value = 490.2857
$ws.Range("M3").Value = -376.2857
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H92").Value = 36316.582
$ws.Range("J92").Value = 36316.582
$ws.Range("L92").Value = 36316.582
$ws.Range("N92").Value = -41308.582
$ws.Range("H96").Value = 25993.889
$ws.Range("I96").Value = 18618.75
$ws.Range("K96").Value = 18618.75
$ws.Range("M96").Value = -15872.75
$ws.Range("H134").Value = 411723.3
$ws.Range("I134").Value = 1893.0186
$ws.Range("J134").Value = 3573271.2
$ws.Range("K134").Value = 5679.0558
$ws.Range("L134").Value = 10719813.6
$ws.Range("M134").Value = -3144.0558
$ws.Range("N134").Value = -10724883.6

# Clear cells that no longer hold a value on BSM
$ws.Range("N40").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Update cell values on CRP
$ws.Range("H16").Value = 3237
$ws.Range("I16").Value = 2566.3333
$ws.Range("J16").Value = 3488.5
$ws.Range("K16").Value = 2566.3333
$ws.Range("L16").Value = 3488.5
$ws.Range("M16").Value = -2279.3333
$ws.Range("N16").Value = -4062.5
$ws.Range("H22").Value = 1120.4166
$ws.Range("I22").Value = 625
$ws.Range("J22").Value = 1368.125
$ws.Range("K22").Value = 625
$ws.Range("L22").Value = 1368.125
$ws.Range("M22").Value = -275
$ws.Range("N22").Value = -2068.125
$ws.Range("H31").Value = 3103.3667
$ws.Range("I31").Value = 1259.25
$ws.Range("J31").Value = 4332.778
$ws.Range("K31").Value = 1259.25
$ws.Range("L31").Value = 4332.778
$ws.Range("M31").Value = -964.25
$ws.Range("N31").Value = -4922.778
$ws.Range("H34").Value = 3103.3667
$ws.Range("I34").Value = 1259.25
$ws.Range("J34").Value = 4332.778
$ws.Range("K34").Value = 1259.25
$ws.Range("L34").Value = 4332.778
$ws.Range("M34").Value = -1057.25
$ws.Range("N34").Value = -4736.778
$ws.Range("H43").Value = 11962.6
$ws.Range("J43").Value = 11962.6
$ws.Range("L43").Value = 11962.6
$ws.Range("N43").Value = -12330.6
$ws.Range("H99").Value = 2480.5386
$ws.Range("I99").Value = 2427.4443
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 2427.4443
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = -929.4443000000001
$ws.Range("N99").Value = -5596
$ws.Range("H101").Value = 11962.6
$ws.Range("J101").Value = 11962.6
$ws.Range("L101").Value = 11962.6
$ws.Range("N101").Value = -18452.6
$ws.Range("H105").Value = 6788.75
$ws.Range("I105").Value = 3685.875
$ws.Range("K105").Value = 3685.875
$ws.Range("M105").Value = -1938.875
$ws.Range("H107").Value = 366.41666
$ws.Range("I107").Value = 308.53845
$ws.Range("J107").Value = 434.81818
$ws.Range("K107").Value = 308.53845
$ws.Range("L107").Value = 434.81818
$ws.Range("M107").Value = 1611.46155
$ws.Range("N107").Value = -4274.81818
$ws.Range("H113").Value = 3237
$ws.Range("I113").Value = 2566.3333
$ws.Range("J113").Value = 3488.5
$ws.Range("K113").Value = 2566.3333
$ws.Range("L113").Value = 3488.5
$ws.Range("M113").Value = -396.3332999999998
$ws.Range("N113").Value = -7828.5
$ws.Range("H126").Value = 2480.5386
$ws.Range("I126").Value = 2427.4443
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 7282.3329
$ws.Range("L126").Value = 7800
$ws.Range("M126").Value = -4812.3329
$ws.Range("N126").Value = -12740
$ws.Range("H134").Value = 1715.8387
$ws.Range("I134").Value = 1492.2593
$ws.Range("K134").Value = 4476.7779
$ws.Range("M134").Value = -1941.7779

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Update cell values on CUL
$ws.Range("H4").Value = 28664832
$ws.Range("I4").Value = 32182798
$ws.Range("J4").Value = 16791704
$ws.Range("K4").Value = 96548394
$ws.Range("L4").Value = 50375112
$ws.Range("M4").Value = -96548282
$ws.Range("N4").Value = -50375336
$ws.Range("H5").Value = 1038.1052
$ws.Range("J5").Value = 2564.3333
$ws.Range("L5").Value = 7692.999899999999
$ws.Range("N5").Value = -7916.999899999999
$ws.Range("H33").Value = 148.81818
$ws.Range("J33").Value = 154
$ws.Range("L33").Value = 924
$ws.Range("N33").Value = -1490
$ws.Range("H39").Value = 2978.9
$ws.Range("I39").Value = 3148.1667
$ws.Range("J39").Value = 2725
$ws.Range("K39").Value = 9444.500100000001
$ws.Range("L39").Value = 8175
$ws.Range("M39").Value = -9150.500100000001
$ws.Range("N39").Value = -8763
$ws.Range("H44").Value = 527.96
$ws.Range("I44").Value = 376.36365
$ws.Range("J44").Value = 647.0714
$ws.Range("K44").Value = 1129.09095
$ws.Range("L44").Value = 1941.2142
$ws.Range("M44").Value = -731.09095
$ws.Range("N44").Value = -2737.2142
$ws.Range("H60").Value = 560.55554
$ws.Range("I60").Value = 89
$ws.Range("J60").Value = 1150
$ws.Range("K60").Value = 267
$ws.Range("L60").Value = 3450
$ws.Range("M60").Value = -16
$ws.Range("N60").Value = -3952
$ws.Range("H63").Value = 2384.1667
$ws.Range("I63").Value = 1768.6666
$ws.Range("K63").Value = 5305.9998
$ws.Range("M63").Value = -4556.9998
$ws.Range("H66").Value = 2384.1667
$ws.Range("I66").Value = 1768.6666
$ws.Range("K66").Value = 15917.9994
$ws.Range("M66").Value = -12173.9994
$ws.Range("H104").Value = 25624.875
$ws.Range("I104").Value = 500
$ws.Range("J104").Value = 40699.8
$ws.Range("K104").Value = 1500
$ws.Range("L104").Value = 122099.4
$ws.Range("M104").Value = 1121
$ws.Range("N104").Value = -127341.4
$ws.Range("H107").Value = 2029.1111
$ws.Range("J107").Value = 1736.2778
$ws.Range("L107").Value = 5208.8334
$ws.Range("N107").Value = -9048.8334
$ws.Range("H113").Value = 1062.6428
$ws.Range("I113").Value = 1694.6
$ws.Range("K113").Value = 5083.799999999999
$ws.Range("M113").Value = -2913.799999999999
$ws.Range("H122").Value = 1084.1
$ws.Range("I122").Value = 656.8
$ws.Range("K122").Value = 5911.2
$ws.Range("M122").Value = -3461.2
$ws.Range("H132").Value = 1748.0385
$ws.Range("J132").Value = 2946
$ws.Range("L132").Value = 26514
$ws.Range("N132").Value = -31574
$ws.Range("H135").Value = 1038.1052
$ws.Range("J135").Value = 2564.3333
$ws.Range("L135").Value = 23078.9997
$ws.Range("N135").Value = -28148.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Update cell values on GSM
$ws.Range("H59").Value = 13999.5
$ws.Range("J59").Value = 13999.5
$ws.Range("L59").Value = 13999.5
$ws.Range("N59").Value = -15165.5
$ws.Range("H80").Value = 66668890
$ws.Range("I80").Value = 71430860
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 71430860
$ws.Range("L80").Value = 1200
$ws.Range("M80").Value = -71429862
$ws.Range("N80").Value = -3196
$ws.Range("H83").Value = 66668890
$ws.Range("I83").Value = 71430860
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 357154300
$ws.Range("L83").Value = 6000
$ws.Range("M83").Value = -357149308
$ws.Range("N83").Value = -15984
$ws.Range("H113").Value = 27786648
$ws.Range("I113").Value = 58831276
$ws.Range("J113").Value = 9873.842000000001
$ws.Range("K113").Value = 58831276
$ws.Range("L113").Value = 9873.842000000001
$ws.Range("M113").Value = -58829106
$ws.Range("N113").Value = -14213.842
$ws.Range("H132").Value = 90911470
$ws.Range("I132").Value = 125002904
$ws.Range("K132").Value = 375008712
$ws.Range("M132").Value = -375006182

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Update cell values on LTW
$ws.Range("H22").Value = 1421.9524
$ws.Range("I22").Value = 617.2353000000001
$ws.Range("J22").Value = 1969.16
$ws.Range("K22").Value = 617.2353000000001
$ws.Range("L22").Value = 1969.16
$ws.Range("M22").Value = -322.2353000000001
$ws.Range("N22").Value = -2559.16
$ws.Range("H27").Value = 1421.9524
$ws.Range("I27").Value = 617.2353000000001
$ws.Range("J27").Value = 1969.16
$ws.Range("K27").Value = 617.2353000000001
$ws.Range("L27").Value = 1969.16
$ws.Range("M27").Value = -510.2353000000001
$ws.Range("N27").Value = -2183.16
$ws.Range("H46").Value = 1101.7778
$ws.Range("I46").Value = 674.6875
$ws.Range("J46").Value = 2153.077
$ws.Range("K46").Value = 674.6875
$ws.Range("L46").Value = 2153.077
$ws.Range("M46").Value = -486.6875
$ws.Range("N46").Value = -2529.077
$ws.Range("H61").Value = 4165.636
$ws.Range("I61").Value = 2392.5454
$ws.Range("J61").Value = 5938.727
$ws.Range("K61").Value = 2392.5454
$ws.Range("L61").Value = 5938.727
$ws.Range("M61").Value = -2190.5454
$ws.Range("N61").Value = -6342.727
$ws.Range("H93").Value = 3681.1086
$ws.Range("I93").Value = 1797.238
$ws.Range("J93").Value = 5263.56
$ws.Range("K93").Value = 1797.238
$ws.Range("L93").Value = 5263.56
$ws.Range("M93").Value = -549.2380000000001
$ws.Range("N93").Value = -7759.56
$ws.Range("H96").Value = 39999
$ws.Range("J96").Value = 39999
$ws.Range("L96").Value = 39999
$ws.Range("N96").Value = -45491
$ws.Range("H113").Value = 4165.636
$ws.Range("I113").Value = 2392.5454
$ws.Range("J113").Value = 5938.727
$ws.Range("K113").Value = 2392.5454
$ws.Range("L113").Value = 5938.727
$ws.Range("M113").Value = -222.5454
$ws.Range("N113").Value = -10278.727
$ws.Range("H122").Value = 3400
$ws.Range("H127").Value = 222000
$ws.Range("J127").Value = 222000
$ws.Range("L127").Value = 222000
$ws.Range("N127").Value = -231920
$ws.Range("H136").Value = 41606.65
$ws.Range("I136").Value = 4093.4285
$ws.Range("K136").Value = 12280.2855
$ws.Range("M136").Value = -9730.2855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Update cell values on WVR
$ws.Range("H24").Value = 15500
$ws.Range("J24").Value = 15500
$ws.Range("L24").Value = 15500
$ws.Range("N24").Value = -15960
$ws.Range("H81").Value = 1079.625
$ws.Range("J81").Value = 1250
$ws.Range("L81").Value = 2500
$ws.Range("N81").Value = -4622
$ws.Range("H84").Value = 1079.625
$ws.Range("J84").Value = 1250
$ws.Range("L84").Value = 12500
$ws.Range("N84").Value = -23108
$ws.Range("H104").Value = 11020.857
$ws.Range("J104").Value = 9524.333000000001
$ws.Range("L104").Value = 9524.333000000001
$ws.Range("N104").Value = -16512.333
$ws.Range("H105").Value = 16500
$ws.Range("J105").Value = 16500
$ws.Range("L105").Value = 16500
$ws.Range("N105").Value = -23488
$ws.Range("H109").Value = 124737.4
$ws.Range("J109").Value = 124737.4
$ws.Range("L109").Value = 124737.4
$ws.Range("N109").Value = -127511.4
$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820
